$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the user (C2) and claim number (E2) on row 2.
# These cells store text that looks numeric, so force text entry with a leading apostrophe.
$ws.Range("C2").Value = "dgariffo"
$ws.Range("E2").Value = "'1120194100385"

# Restore the selection to E2 (was E4).
$ws.Range("E2").Select()
